# Update "想去人数" (people wanting to go) counts on the "展览" (Exhibition)
# and "全部类型" (All Types) sheets to match the latest scrape.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (sheet index 1) ----
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F2").Value  = 0      # was 318
$wsExpo.Range("F4").Value  = 0      # was 8319
$wsExpo.Range("F5").Value  = 6068   # was 6066
$wsExpo.Range("F6").Value  = 0      # was 522
$wsExpo.Range("F9").Value  = 0      # was 71
$wsExpo.Range("F11").Value = 0      # was 1012
$wsExpo.Range("F12").Value = 0      # was 81

# ---- Sheet "全部类型" (sheet index 4) ----
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value  = 0      # was 318
$wsAll.Range("F3").Value  = 0      # was 20
$wsAll.Range("F4").Value  = 0      # was 8319
$wsAll.Range("F5").Value  = 6068   # was 6066
$wsAll.Range("F6").Value  = 0      # was 522
$wsAll.Range("F8").Value  = 0      # was 19
$wsAll.Range("F9").Value  = 0      # was 71
$wsAll.Range("F10").Value = 0      # was 315
$wsAll.Range("F11").Value = 0      # was 96
$wsAll.Range("F12").Value = 0      # was 1
$wsAll.Range("F14").Value = 0      # was 7
$wsAll.Range("F15").Value = 1019   # was 1012
$wsAll.Range("F16").Value = 0      # was 81
$wsAll.Range("F17").Value = 0      # was 2
